$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (co2)
$ws.Range("C2").Value = 552.4246960852229
$ws.Range("D2").Value = 136.2986616245581
$ws.Range("F2").Value = 449
$ws.Range("G2").Value = 509
$ws.Range("H2").Value = 622

# Row 3 (humidity)
$ws.Range("C3").Value = 40.22822938490521
$ws.Range("D3").Value = 5.303746036525993
$ws.Range("F3").Value = 37.12
$ws.Range("G3").Value = 39.94
$ws.Range("H3").Value = 43.28

# Row 4 (pm25)
$ws.Range("C4").Value = 1.599382867147636
$ws.Range("D4").Value = 2.134985795785292
$ws.Range("F4").Value = 0.55
$ws.Range("G4").Value = 1.06
$ws.Range("H4").Value = 2.02

# Row 5 (pressure)
$ws.Range("C5").Value = 322.8763409526746
$ws.Range("D5").Value = 11.37829707917052
$ws.Range("F5").Value = 315.06
$ws.Range("G5").Value = 324.32
$ws.Range("H5").Value = 332.05

# Row 6 (temperature)
$ws.Range("C6").Value = 20.8051068664887
$ws.Range("D6").Value = 2.521383329803396
$ws.Range("F6").Value = 19.44
$ws.Range("H6").Value = 22.27

# Row 7 (rssi)
$ws.Range("C7").Value = -76.07091309600243
$ws.Range("D7").Value = 22.90160120971134
$ws.Range("I7").Value = -28

# Row 8 (snr)
$ws.Range("C8").Value = 7.637385816818627
$ws.Range("D8").Value = 6.897731999772141

# Row 9 (SF)
$ws.Range("C9").Value = 9.321617443314297
$ws.Range("D9").Value = 1.685160368922639

# Row 10 (frequency)
$ws.Range("C10").Value = 867.8302478626704
$ws.Range("D10").Value = 0.4614263994982076

# Row 11 (toa)
$ws.Range("C11").Value = 0.555451883080458
$ws.Range("D11").Value = 0.5887337633082016

# Row 12 (distance)
$ws.Range("C12").Value = 22.7386924272632
$ws.Range("D12").Value = 12.29186826816025

# Row 13 (c_walls)
$ws.Range("C13").Value = 0.673922042374886
$ws.Range("D13").Value = 0.7505261244078241

# Row 14 (w_walls)
$ws.Range("C14").Value = 1.826762241070524
$ws.Range("D14").Value = 1.664037666823125

# Row 15 (exp_pl)
$ws.Range("C15").Value = 93.47091309600226
$ws.Range("D15").Value = 22.90160120971134
$ws.Range("E15").Value = 45.4

# Row 16 (n_power)
$ws.Range("C16").Value = -85.30844803363917
$ws.Range("D16").Value = 20.60786016223266
$ws.Range("H16").Value = -67.8707776445072

# Row 17 (esp)
$ws.Range("C17").Value = -77.67106221682054
$ws.Range("D17").Value = 25.31265400822541
$ws.Range("G17").Value = -72.75746206410165
$ws.Range("H17").Value = -57.25410721860875
$ws.Range("I17").Value = -28.39612087980607
